$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that Excel would otherwise auto-convert (numbers, percents)
# as literal text, preserving the original (default) cell style.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "306.31"
Set-TextValue $ws.Range("E2") "-0.69%"
Set-TextValue $ws.Range("D3") "38.89"
Set-TextValue $ws.Range("E3") "7.08%"
Set-TextValue $ws.Range("D4") "5.093"
Set-TextValue $ws.Range("E4") "0.61%"
Set-TextValue $ws.Range("D5") "0.08051"
Set-TextValue $ws.Range("E5") "-0.88%"
Set-TextValue $ws.Range("D6") "1.941"
Set-TextValue $ws.Range("E6") "-7.29%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D7") "4.195"
Set-TextValue $ws.Range("E7") "0.94%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws.Range("D8") "8.001"
Set-TextValue $ws.Range("E8") "1.90%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D9") "0.9315"
Set-TextValue $ws.Range("E9") "0.12%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.1442"
Set-TextValue $ws.Range("E10") "0.03%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1926"
Set-TextValue $ws.Range("E11") "-0.35%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.09061"
Set-TextValue $ws.Range("E12") "-0.03%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03513"
Set-TextValue $ws.Range("E13") "1.61%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09770"
Set-TextValue $ws.Range("E14") "-1.37%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001399"
Set-TextValue $ws.Range("E15") "-0.63%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.005977"
Set-TextValue $ws.Range("E16") "-5.36%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.787"
Set-TextValue $ws.Range("E17") "-1.30%"
Set-TextValue $ws.Range("D18") "3.434"
Set-TextValue $ws.Range("E18") "1.71%"
Set-TextValue $ws.Range("E19") "-0.72%"
Set-TextValue $ws.Range("E20") "1.56%"
Set-TextValue $ws.Range("D21") "4.892"
Set-TextValue $ws.Range("E21") "1.91%"
Set-TextValue $ws.Range("D22") "0.2415"
Set-TextValue $ws.Range("E22") "3.17%"
Set-TextValue $ws.Range("E23") "0.48%"
Set-TextValue $ws.Range("D24") "0.001237"
Set-TextValue $ws.Range("E24") "0.24%"
Set-TextValue $ws.Range("D26") "0.0001302"
Set-TextValue $ws.Range("D39") "0.02033"
Set-TextValue $ws.Range("E39") "0.69%"
Set-TextValue $ws.Range("D40") "0.05048"
Set-TextValue $ws.Range("E40") "-2.31%"
Set-TextValue $ws.Range("D41") "0.007429"
Set-TextValue $ws.Range("E41") "-0.88%"
Set-TextValue $ws.Range("D42") "0.01017"
Set-TextValue $ws.Range("E42") "0.41%"
Set-TextValue $ws.Range("D43") "0.1348"
Set-TextValue $ws.Range("E43") "-1.56%"
Set-TextValue $ws.Range("D44") "0.002123"
Set-TextValue $ws.Range("E44") "-2.17%"
Set-TextValue $ws.Range("D45") "0.009119"
Set-TextValue $ws.Range("D46") "0.00006193"
Set-TextValue $ws.Range("E46") "-1.38%"
Set-TextValue $ws.Range("E47") "0.14%"
Set-TextValue $ws.Range("D48") "0.003075"
Set-TextValue $ws.Range("E49") "28.12%"
Set-TextValue $ws.Range("D50") "0.00002103"
Set-TextValue $ws.Range("E50") "0.14%"
Set-TextValue $ws.Range("D51") "0.0002003"
Set-TextValue $ws.Range("E51") "0.14%"
